$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TRANSIENT")

# The simulation name formula in E3 now appends a "case3" suffix via TEXTJOIN,
# turning e.g. "TEND_20000_STPMIN_10_NELEMS_200" into
# "TEND_20000_STPMIN_10_NELEMS_200_case3".
$ws.Range("E3").Formula = '=_xlfn.TEXTJOIN("_",TRUE,A6,E6,A8,E8,[1]GRID!$A$4,[1]GRID!$E$4,"case3")'
